# Working on regression of data sets
# Add a new "sample size" column (M) to the Ducry1979 / Thompson1965 data table
# and bold the chart axis titles + shrink the scatter markers slightly on both
# liver-weight and body-weight charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M: "sample size" / "n" ---------------------------------
$ws.Range("M10").Value = "sample size"
$ws.Range("M10").Style = $ws.Range("L10").Style

$ws.Range("M11").Value = "n"
$ws.Range("M11").Style = $ws.Range("L11").Style

$ws.Range("M12").Value = 25
$ws.Range("M13").Value = 25
$ws.Range("M14").Value = 25
$ws.Range("M15").Value = 43
$ws.Range("M16").Value = 64
$ws.Range("M17").Value = 53
$ws.Range("M18").Value = 21

# Leave the selection where the user ended up after typing the data in.
$ws.Range("A1").Select()
$ws.Range("M18").Select()

# --- Chart formatting tweaks --------------------------------------------
foreach ($co in $ws.ChartObjects()) {
    $chart = $co.Chart
    foreach ($axis in @($chart.Axes(1), $chart.Axes(2))) {
        $axis.AxisTitle.Font.Bold = $true
    }
    foreach ($series in $chart.SeriesCollection()) {
        $series.MarkerSize = 6
    }
}
